# "upgrade left table until javakheti"
# Rename sheet, blank out older-year Urban/Rural figures (replacing them
# with the confidentiality marker already used elsewhere on the sheet),
# normalize the ellipsis marker text, and remove the stray blank row
# that used to separate the data table from the footnote.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet tab from "1" to "Gurjaani"
$ws.Name = "Gurjaani"

# 2) Normalize the confidentiality marker glyph "…" -> "..." everywhere
#    it is used on the sheet (in place, so the shared string itself is
#    updated rather than creating a duplicate string entry).
$ws.Cells.Replace("…", "...") | Out-Null

# Grab the already-normalized marker value (used e.g. by N5) so that the
# newly-blanked cells reuse the very same string value.
$marker = $ws.Cells.Item(5, 14).Value()

# 3) Row 6 ("Urban"): blank out 2010-2020 (cols B-L), keep 2021 (col M)
for ($c = 2; $c -le 12; $c++) {
    $ws.Cells.Item(6, $c).Value = $marker
}

# 4) Row 7 ("Rural"): blank out 2010-2015 & 2017-2020 (cols B-G, I-L),
#    keep 2016 (col H) and 2021 (col M)
for ($c = 2; $c -le 7; $c++) {
    $ws.Cells.Item(7, $c).Value = $marker
}
for ($c = 9; $c -le 12; $c++) {
    $ws.Cells.Item(7, $c).Value = $marker
}

# 5) Remove the stray blank row 8 so the footnote row (old row 9) moves
#    up to row 8, matching the tightened table.
$ws.Rows.Item(8).Delete()
